$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new values for columns C (nombre_aides) and E (montant_total)
$updates = @(
    @{ Row = 13;  C = 187870;  E = 1168643288 },
    @{ Row = 91;  C = 18887;   E = 75398735 },
    @{ Row = 115; C = 81811;   E = 436758478 },
    @{ Row = 121; C = 1306458; E = 2275688370 },
    @{ Row = 127; C = 9163;    E = 110743559 },
    @{ Row = 129; C = 633923;  E = 3437003883 },
    @{ Row = 132; C = 586118;  E = 3475155346 },
    @{ Row = 134; C = 7028;    E = 16907539 },
    @{ Row = 136; C = 26709;   E = 144431319 },
    @{ Row = 161; C = 555;     E = 1818521 },
    @{ Row = 171; C = 95830;   E = 490707985 },
    @{ Row = 174; C = 40453;   E = 240012967 },
    @{ Row = 186; C = 236844;  E = 1190209648 },
    @{ Row = 215; C = 230266;  E = 408781748 }
)

foreach ($u in $updates) {
    $ws.Range("C" + $u.Row).Value = $u.C
    $ws.Range("E" + $u.Row).Value = $u.E
}

$wb.Save()
